$d = $word.ActiveDocument

# --- Build the colored/highlighted "<---" + message runs in an isolated
# scratch paragraph appended at the end of the document. We do this in
# isolation because this runtime's Range.HighlightColorIndex setter applies
# to the *entire paragraph* that contains the range rather than just the
# range itself; building it in its own paragraph (where the paragraph's
# full extent *is* the text we want highlighted) sidesteps that limitation.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)

$arrowStart = $end.Start
$end.InsertAfter("<---")
$arrowEnd = $end.End
$end.Collapse(0)
$end.InsertAfter("M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0")

# Force a run split between the "<---" run and the message run (they must
# stay distinct <w:r> elements even though their formatting ends up
# identical) by toggling a property on just the first part and back off.
$arrowRng = $d.Range($arrowStart, $arrowEnd)
$arrowRng.Bold = 1
$arrowRng.Bold = 0

# Apply the target character formatting. At this point the scratch
# paragraph's range is *exactly* the "<---" + message text, so the
# paragraph-wide HighlightColorIndex quirk ends up scoped correctly.
$scratchPara = $d.Paragraphs($d.Paragraphs.Count).Range
$scratchPara.HighlightColorIndex = 16
$scratchPara.Font.Color = 42495
$scratchPara.Font.Size = 16

# Surround with plain (unformatted) 4-space runs.
$scratchPara2 = $d.Paragraphs($d.Paragraphs.Count).Range
$scratchPara2.InsertBefore("    ")
$scratchPara3 = $d.Paragraphs($d.Paragraphs.Count).Range
$scratchPara3.InsertAfter("    ")

# Exclude the trailing paragraph mark from the range we copy.
$scratchPara3.MoveEnd(1, -1)
$payload = $scratchPara3.FormattedText

# --- Insert the payload right after "Basic " (before "if") in paragraph 1.
$target = $d.Content
$target.Find.Execute("Basic ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.FormattedText = $payload

# --- Clean up: remove the scratch paragraph (including its paragraph mark)
# we used to build/format the payload.
$scratchPara4 = $d.Paragraphs($d.Paragraphs.Count).Range
$scratchPara4.Delete()
$trailingMark = $d.Range($scratchPara4.Start - 1, $scratchPara4.Start)
$trailingMark.Delete()
